$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency Price (D) and Volume(1h) (E) columns with refreshed data.
# Cells hold text values (e.g. "328.07", "-0.08%"), so force Text format first
# to stop Excel auto-converting the assignment into a number/percentage.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "328.07"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.08%"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "44.01"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-0.02%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.382"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-2.29%"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08381"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "4.08%"

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-6.13%"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9766"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "2.32%"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1122"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "0.04%"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1900"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "1.26%"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09785"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-1.40%"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.04614"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-2.45%"

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.54%"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001291"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "1.84%"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.006147"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.58%"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.404"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.91%"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.434"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.28%"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.3330"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.72%"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.162"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-9.92%"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1370"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-2.31%"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.2550"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "2.13%"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.04156"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "1.47%"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.001296"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-1.13%"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.004410"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "1.83%"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0001301"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "3.83%"

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02653"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "1.37%"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05641"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "0.63%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.007811"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "2.54%"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1415"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "1.27%"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.007364"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-0.19%"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002112"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "6.11%"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007913"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-10.64%"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3514"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006841"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-3.82%"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000751"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.16%"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003514"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "2.31%"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.003534"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.98%"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002102"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.16%"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002002"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.16%"
